$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 48, shifting existing rows 48+ down by one
$ws.Rows.Item(48).Insert()

# Populate the new row 48 with the new record
$ws.Cells.Item(48, 1).Value = 6
$ws.Cells.Item(48, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(48, 3).Value = "Metropolitana"
$ws.Cells.Item(48, 4).Value = 45082
$ws.Cells.Item(48, 5).Value = 13
$ws.Cells.Item(48, 6).Value = 100114007
$ws.Cells.Item(48, 7).Value = "Jengibre"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 280
$ws.Cells.Item(48, 11).Value = 15000
$ws.Cells.Item(48, 12).Value = 16000
$ws.Cells.Item(48, 13).Value = 15571
$ws.Cells.Item(48, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(48, 15).Value = "Perú"
$ws.Cells.Item(48, 16).Value = 1198
$ws.Cells.Item(48, 17).Value = 13
$ws.Cells.Item(48, 18).Value = "Hortaliza"
